# Helper: VBA-style RGB() packs bytes as 0x00BBGGRR (the native encoding
# used by the PowerPoint object model's RGB color properties).
function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# -----------------------------------------------------------------------
# 1) Swap the presentation's active theme palette ("Integral" / Red Violet)
#    for the Office default palette (the palette previously only used by
#    the notes master theme). The font scheme and format scheme are
#    already identical between the two themes, so only the 12 theme
#    colors (and anything using scheme colors) need to change.
# -----------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

$officePalette = @(
    @(0x00, 0x00, 0x00),  # dk1
    @(0xFF, 0xFF, 0xFF),  # lt1
    @(0x44, 0x54, 0x6A),  # dk2
    @(0xE7, 0xE6, 0xE6),  # lt2
    @(0x5B, 0x9B, 0xD5),  # accent1
    @(0xED, 0x7D, 0x31),  # accent2
    @(0xA5, 0xA5, 0xA5),  # accent3
    @(0xFF, 0xC0, 0x00),  # accent4
    @(0x44, 0x72, 0xC4),  # accent5
    @(0x70, 0xAD, 0x47),  # accent6
    @(0x05, 0x63, 0xC1),  # hlink
    @(0x95, 0x4F, 0x72)   # folHlink
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $rgb = $officePalette[$i - 1]
    $themeColors.Colors($i).RGB = RGBVal $rgb[0] $rgb[1] $rgb[2]
}

# -----------------------------------------------------------------------
# 2) Change the table style used by the table on slide 5 to the other
#    built-in table style.
# -----------------------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{AF42E18C-8202-43C6-8E03-75CC8C1CF539}")
